$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 73282.86
$ws.Range("I113").Value = 334734.66
$ws.Range("J113").Value = 1977.8182
$ws.Range("K113").Value = 334734.66
$ws.Range("L113").Value = 1977.8182
$ws.Range("M113").Value = -331480.66
$ws.Range("N113").Value = -8485.8182

$ws.Range("H115").Value = 20259
$ws.Range("I115").Value = 20259
$ws.Range("K115").Value = 60777
$ws.Range("M115").Value = -59210

$ws.Range("H125").Value = 1877.2222
$ws.Range("I125").Value = 3821.3333
$ws.Range("J125").Value = 1488.4
$ws.Range("K125").Value = 34391.9997
$ws.Range("L125").Value = 13395.6
$ws.Range("M125").Value = -31931.9997
$ws.Range("N125").Value = -18315.6

$ws.Range("H127").Value = 1588.3889
$ws.Range("I127").Value = 372.2
$ws.Range("J127").Value = 2056.1538
$ws.Range("K127").Value = 1116.6
$ws.Range("L127").Value = 6168.4614
$ws.Range("M127").Value = 3843.4
$ws.Range("N127").Value = -16088.4614

$ws.Range("H129").Value = 3651.0557
$ws.Range("J129").Value = 1082.7391
$ws.Range("L129").Value = 3248.2173
$ws.Range("N129").Value = -13248.2173

$ws.Range("H137").Value = 1636.2059
$ws.Range("I137").Value = 1334.4
$ws.Range("J137").Value = 3899.75
$ws.Range("K137").Value = 4003.2
$ws.Range("L137").Value = 11699.25
$ws.Range("M137").Value = -1453.2
$ws.Range("N137").Value = -16799.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 113352.445
$ws.Range("I45").Value = 144030.28
$ws.Range("K45").Value = 144030.28
$ws.Range("M45").Value = -143653.28

$ws.Range("H74").Value = 931.2308
$ws.Range("I74").Value = 991.7778
$ws.Range("J74").Value = 795
$ws.Range("K74").Value = 991.7778
$ws.Range("L74").Value = 795
$ws.Range("M74").Value = -117.7778
$ws.Range("N74").Value = -2543

$ws.Range("H77").Value = 931.2308
$ws.Range("I77").Value = 991.7778
$ws.Range("J77").Value = 795
$ws.Range("K77").Value = 4958.889
$ws.Range("L77").Value = 3975
$ws.Range("M77").Value = -590.8890000000001
$ws.Range("N77").Value = -12711

$ws.Range("H101").Value = 35695
$ws.Range("J101").Value = 35695
$ws.Range("L101").Value = 35695
$ws.Range("N101").Value = -42185

$ws.Range("H132").Value = 2045.9149
$ws.Range("I132").Value = 1610
$ws.Range("J132").Value = 3317.3333
$ws.Range("K132").Value = 4830
$ws.Range("L132").Value = 9951.999899999999
$ws.Range("M132").Value = -2300
$ws.Range("N132").Value = -15011.9999

$ws.Range("H141").Value = 48371.6
$ws.Range("J141").Value = 48371.6
$ws.Range("L141").Value = 48371.6
$ws.Range("N141").Value = -58731.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 377.66666
$ws.Range("I22").Value = 344.25
$ws.Range("K22").Value = 344.25
$ws.Range("M22").Value = -171.25

$ws.Range("H86").Value = 45301.73
$ws.Range("I86").Value = 55241.19
$ws.Range("J86").Value = 3556
$ws.Range("K86").Value = 55241.19
$ws.Range("L86").Value = 3556
$ws.Range("M86").Value = -54118.19
$ws.Range("N86").Value = -5802

$ws.Range("H89").Value = 45301.73
$ws.Range("I89").Value = 55241.19
$ws.Range("J89").Value = 3556
$ws.Range("K89").Value = 276205.95
$ws.Range("L89").Value = 17780
$ws.Range("M89").Value = -270589.95
$ws.Range("N89").Value = -29012

$ws.Range("H107").Value = 43496024
$ws.Range("I107").Value = 47638420
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 47638420
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = -47636500
$ws.Range("N107").Value = -4740

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1214.1428
$ws.Range("J16").Value = 1166.3334
$ws.Range("L16").Value = 1166.3334
$ws.Range("N16").Value = -1740.3334

$ws.Range("H31").Value = 30276.22
$ws.Range("I31").Value = 1331.1052
$ws.Range("J31").Value = 48016.773
$ws.Range("K31").Value = 1331.1052
$ws.Range("L31").Value = 48016.773
$ws.Range("M31").Value = -1036.1052
$ws.Range("N31").Value = -48606.773

$ws.Range("H34").Value = 30276.22
$ws.Range("I34").Value = 1331.1052
$ws.Range("J34").Value = 48016.773
$ws.Range("K34").Value = 1331.1052
$ws.Range("L34").Value = 48016.773
$ws.Range("M34").Value = -1129.1052
$ws.Range("N34").Value = -48420.773

$ws.Range("H47").Value = 34849.5
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 34849.5
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 34849.5
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -35981.5

$ws.Range("H58").Value = 989.5814
$ws.Range("I58").Value = 989.5814
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 989.5814
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -786.5814
$ws.Range("N58").ClearContents()

$ws.Range("H86").Value = 2936.5715
$ws.Range("I86").Value = 2899.25
$ws.Range("K86").Value = 2899.25
$ws.Range("M86").Value = -1776.25

$ws.Range("H89").Value = 2936.5715
$ws.Range("I89").Value = 2899.25
$ws.Range("K89").Value = 14496.25
$ws.Range("M89").Value = -8880.25

$ws.Range("H105").Value = 1079.05
$ws.Range("I105").Value = 1066.1875
$ws.Range("J105").Value = 1130.5
$ws.Range("K105").Value = 1066.1875
$ws.Range("L105").Value = 1130.5
$ws.Range("M105").Value = 680.8125
$ws.Range("N105").Value = -4624.5

$ws.Range("H106").Value = 22777.143
$ws.Range("J106").Value = 22777.143
$ws.Range("L106").Value = 22777.143
$ws.Range("N106").Value = -25301.143

$ws.Range("H107").Value = 747.3
$ws.Range("I107").Value = 850.38464
$ws.Range("J107").Value = 555.8570999999999
$ws.Range("K107").Value = 850.38464
$ws.Range("L107").Value = 555.8570999999999
$ws.Range("M107").Value = 1069.61536
$ws.Range("N107").Value = -4395.8571

$ws.Range("H110").Value = 39999
$ws.Range("J110").Value = 39999
$ws.Range("L110").Value = 39999
$ws.Range("N110").Value = -48179

$ws.Range("H113").Value = 1214.1428
$ws.Range("J113").Value = 1166.3334
$ws.Range("L113").Value = 1166.3334
$ws.Range("N113").Value = -5506.3334

$ws.Range("H136").Value = 989.5814
$ws.Range("I136").Value = 989.5814
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2968.7442
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -418.7442000000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1840.3334
$ws.Range("I121").Value = 1619.4
$ws.Range("J121").Value = 2116.5
$ws.Range("K121").Value = 4858.200000000001
$ws.Range("L121").Value = 6349.5
$ws.Range("M121").Value = -3548.200000000001
$ws.Range("N121").Value = -8969.5

$ws.Range("H127").Value = 1060
$ws.Range("J127").Value = 1060
$ws.Range("L127").Value = 3180
$ws.Range("N127").Value = -13100

$ws.Range("H131").Value = 7319.387
$ws.Range("J131").Value = 7826.0815
$ws.Range("L131").Value = 23478.2445
$ws.Range("N131").Value = -33558.2445

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 116128.78
$ws.Range("J70").Value = 6244.8887
$ws.Range("L70").Value = 6244.8887
$ws.Range("N70").Value = -6784.8887

$ws.Range("H73").Value = 116128.78
$ws.Range("J73").Value = 6244.8887
$ws.Range("L73").Value = 6244.8887
$ws.Range("N73").Value = -8116.8887

$ws.Range("H102").Value = 2526.7222
$ws.Range("I102").Value = 1762.875
$ws.Range("K102").Value = 1762.875
$ws.Range("M102").Value = -140.875

$ws.Range("H132").Value = 3283.652
$ws.Range("I132").Value = 2489.5334
$ws.Range("J132").Value = 4772.625
$ws.Range("K132").Value = 7468.600199999999
$ws.Range("L132").Value = 14317.875
$ws.Range("M132").Value = -4938.600199999999
$ws.Range("N132").Value = -19377.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 16894.555
$ws.Range("J104").Value = 16894.555
$ws.Range("L104").Value = 16894.555
$ws.Range("N104").Value = -23882.555

$ws.Range("H125").Value = 20000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 20000
$ws.Range("N125").Value = -29840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2720149.2
$ws.Range("J62").Value = 2726.923
$ws.Range("L62").Value = 2726.923
$ws.Range("N62").Value = -3974.923

$ws.Range("H65").Value = 2720149.2
$ws.Range("J65").Value = 2726.923
$ws.Range("L65").Value = 13634.615
$ws.Range("N65").Value = -19874.615

$ws.Range("H107").Value = 125571.625
$ws.Range("I107").Value = 599.6667
$ws.Range("J107").Value = 200554.8
$ws.Range("K107").Value = 1799.0001
$ws.Range("L107").Value = 601664.3999999999
$ws.Range("M107").Value = 120.9999
$ws.Range("N107").Value = -605504.3999999999
